$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
  # hunk
  $ws.Cells.Item(18, 8).Value = 166667180
  $ws.Cells.Item(18, 10).Value = 500001000
  $ws.Cells.Item(18, 12).Value = 500001000
  $ws.Cells.Item(18, 14).Value = -500001568
  # hunk
  $ws.Cells.Item(40, 8).Value = 3960.1667
  $ws.Cells.Item(40, 9).Value = 2598.3333
  $ws.Cells.Item(40, 11).Value = 2598.3333
  $ws.Cells.Item(40, 13).Value = -2423.3333
  # hunk
  $ws.Cells.Item(43, 8).Value = 0
  $ws.Cells.Item(43, 10).Value = 0
  $ws.Cells.Item(43, 12).Value = 0
  $ws.Cells.Item(43, 14).ClearContents()
  # hunk
  $ws.Cells.Item(107, 8).Value = 814.8333
  $ws.Cells.Item(107, 9).Value = 814.8333
  $ws.Cells.Item(107, 11).Value = 814.8333
  $ws.Cells.Item(107, 13).Value = 1105.1667
  # hunk
  $ws.Cells.Item(138, 8).Value = 3433.1294
  $ws.Cells.Item(138, 9).Value = 1889.4828
  $ws.Cells.Item(138, 10).Value = 4232.518
  $ws.Cells.Item(138, 11).Value = 5668.4484
  $ws.Cells.Item(138, 12).Value = 12697.554
  $ws.Cells.Item(138, 13).Value = -528.4484000000002
  $ws.Cells.Item(138, 14).Value = -22977.554

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
  # hunk
  $ws.Cells.Item(2, 8).Value = 3319.5518
  $ws.Cells.Item(2, 9).Value = 2171.3157
  $ws.Cells.Item(2, 11).Value = 2171.3157
  $ws.Cells.Item(2, 13).Value = -2058.3157
  # hunk
  $ws.Cells.Item(5, 8).Value = 71.28570999999999
  $ws.Cells.Item(5, 9).Value = 71.28570999999999
  $ws.Cells.Item(5, 11).Value = 71.28570999999999
  $ws.Cells.Item(5, 13).Value = 40.71429000000001
  # hunk
  $ws.Cells.Item(32, 8).Value = 1212328.6
  $ws.Cells.Item(32, 9).Value = 572905.6
  $ws.Cells.Item(32, 11).Value = 572905.6
  $ws.Cells.Item(32, 13).Value = -572618.6
  # hunk
  $ws.Cells.Item(74, 8).Value = 21281246
  $ws.Cells.Item(74, 9).Value = 2678.4
  $ws.Cells.Item(74, 11).Value = 2678.4
  $ws.Cells.Item(74, 13).Value = -1804.4
  # hunk
  $ws.Cells.Item(77, 8).Value = 21281246
  $ws.Cells.Item(77, 9).Value = 2678.4
  $ws.Cells.Item(77, 11).Value = 13392
  $ws.Cells.Item(77, 13).Value = -9024
  # hunk
  $ws.Cells.Item(116, 8).Value = 3319.5518
  $ws.Cells.Item(116, 9).Value = 2171.3157
  $ws.Cells.Item(116, 11).Value = 2171.3157
  $ws.Cells.Item(116, 13).Value = 122.6842999999999
  # hunk
  $ws.Cells.Item(122, 8).Value = 3434.7778
  $ws.Cells.Item(122, 9).Value = 2891.111
  $ws.Cells.Item(122, 11).Value = 8673.332999999999
  $ws.Cells.Item(122, 13).Value = -6223.332999999999
  # hunk
  $ws.Cells.Item(132, 8).Value = 3080979.2
  $ws.Cells.Item(132, 9).Value = 3666118.5
  $ws.Cells.Item(132, 10).Value = 8999.25
  $ws.Cells.Item(132, 11).Value = 10998355.5
  $ws.Cells.Item(132, 12).Value = 26997.75
  $ws.Cells.Item(132, 13).Value = -10995825.5
  $ws.Cells.Item(132, 14).Value = -32057.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
  # hunk
  $ws.Cells.Item(3, 8).Value = 3319.5518
  $ws.Cells.Item(3, 9).Value = 2171.3157
  $ws.Cells.Item(3, 11).Value = 2171.3157
  $ws.Cells.Item(3, 13).Value = -2057.3157
  # hunk
  $ws.Cells.Item(4, 8).Value = 71.28570999999999
  $ws.Cells.Item(4, 9).Value = 71.28570999999999
  $ws.Cells.Item(4, 11).Value = 71.28570999999999
  $ws.Cells.Item(4, 13).Value = 43.71429000000001
  # hunk
  $ws.Cells.Item(22, 8).Value = 649
  $ws.Cells.Item(22, 9).Value = 649
  $ws.Cells.Item(22, 11).Value = 649
  $ws.Cells.Item(22, 13).Value = -476
  # hunk
  $ws.Cells.Item(105, 9).Value = 500975.7
  $ws.Cells.Item(105, 11).Value = 500975.7
  $ws.Cells.Item(105, 13).Value = -499228.7
  # hunk
  $ws.Cells.Item(134, 8).Value = 1969.4445
  $ws.Cells.Item(134, 9).Value = 1501.5161
  $ws.Cells.Item(134, 11).Value = 4504.5483
  $ws.Cells.Item(134, 13).Value = -1969.5483

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
  # hunk
  $ws.Cells.Item(22, 8).Value = 3331.3333
  $ws.Cells.Item(22, 10).Value = 0
  $ws.Cells.Item(22, 12).Value = 0
  $ws.Cells.Item(22, 14).ClearContents()
  # hunk
  $ws.Cells.Item(31, 8).Value = 2406265.5
  $ws.Cells.Item(31, 9).Value = 1377.9584
  $ws.Cells.Item(31, 11).Value = 1377.9584
  $ws.Cells.Item(31, 13).Value = -1082.9584
  # hunk
  $ws.Cells.Item(34, 8).Value = 2406265.5
  $ws.Cells.Item(34, 9).Value = 1377.9584
  $ws.Cells.Item(34, 11).Value = 1377.9584
  $ws.Cells.Item(34, 13).Value = -1175.9584
  # hunk
  $ws.Cells.Item(105, 8).Value = 2343.8333
  $ws.Cells.Item(105, 9).Value = 1965.875
  $ws.Cells.Item(105, 10).Value = 3099.75
  $ws.Cells.Item(105, 11).Value = 1965.875
  $ws.Cells.Item(105, 12).Value = 3099.75
  $ws.Cells.Item(105, 13).Value = -218.875
  $ws.Cells.Item(105, 14).Value = -6593.75
  # hunk
  $ws.Cells.Item(107, 8).Value = 623.5454999999999
  $ws.Cells.Item(107, 9).Value = 387
  $ws.Cells.Item(107, 10).Value = 1037.5
  $ws.Cells.Item(107, 11).Value = 387
  $ws.Cells.Item(107, 12).Value = 1037.5
  $ws.Cells.Item(107, 13).Value = 1533
  $ws.Cells.Item(107, 14).Value = -4877.5
  # hunk
  $ws.Cells.Item(132, 8).Value = 2359.5952
  $ws.Cells.Item(132, 9).Value = 2423.7896
  $ws.Cells.Item(132, 10).Value = 1749.75
  $ws.Cells.Item(132, 11).Value = 7271.3688
  $ws.Cells.Item(132, 12).Value = 5249.25
  $ws.Cells.Item(132, 13).Value = -4741.3688
  $ws.Cells.Item(132, 14).Value = -10309.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
  # hunk
  $ws.Cells.Item(8, 8).Value = 5354.3335
  $ws.Cells.Item(8, 9).Value = 5354.3335
  $ws.Cells.Item(8, 11).Value = 16063.0005
  $ws.Cells.Item(8, 13).Value = -15924.0005
  # hunk
  $ws.Cells.Item(12, 8).Value = 27.583334
  $ws.Cells.Item(12, 9).Value = 7.3333335
  $ws.Cells.Item(12, 11).Value = 22.0000005
  $ws.Cells.Item(12, 13).Value = 150.9999995
  # hunk
  $ws.Cells.Item(68, 8).Value = 1926258.8
  $ws.Cells.Item(68, 9).Value = 1907.875
  $ws.Cells.Item(68, 10).Value = 2276140.8
  $ws.Cells.Item(68, 11).Value = 5723.625
  $ws.Cells.Item(68, 12).Value = 6828422.399999999
  $ws.Cells.Item(68, 13).Value = -4912.625
  $ws.Cells.Item(68, 14).Value = -6830044.399999999
  # hunk
  $ws.Cells.Item(71, 8).Value = 1926258.8
  $ws.Cells.Item(71, 9).Value = 1907.875
  $ws.Cells.Item(71, 10).Value = 2276140.8
  $ws.Cells.Item(71, 11).Value = 17170.875
  $ws.Cells.Item(71, 12).Value = 20485267.2
  $ws.Cells.Item(71, 13).Value = -13114.875
  $ws.Cells.Item(71, 14).Value = -20493379.2
  # hunk
  $ws.Cells.Item(120, 8).Value = 10666
  $ws.Cells.Item(120, 9).Value = 10666
  $ws.Cells.Item(120, 11).Value = 31998
  $ws.Cells.Item(120, 13).Value = -27160
  # hunk
  $ws.Cells.Item(122, 8).Value = 1416.8462
  $ws.Cells.Item(122, 10).Value = 2078.875
  $ws.Cells.Item(122, 12).Value = 18709.875
  $ws.Cells.Item(122, 14).Value = -23609.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
  # hunk
  $ws.Cells.Item(97, 8).Value = 783.26086
  $ws.Cells.Item(97, 9).Value = 572.1905
  $ws.Cells.Item(97, 11).Value = 572.1905
  $ws.Cells.Item(97, 13).Value = -76.19050000000004
  # hunk
  $ws.Cells.Item(122, 8).Value = 71435570
  $ws.Cells.Item(122, 9).Value = 6495.875
  $ws.Cells.Item(122, 11).Value = 19487.625
  $ws.Cells.Item(122, 13).Value = -17037.625
  # hunk
  $ws.Cells.Item(132, 8).Value = 1682.2
  $ws.Cells.Item(132, 9).Value = 1658.8
  $ws.Cells.Item(132, 11).Value = 4976.4
  $ws.Cells.Item(132, 13).Value = -2446.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
  # hunk
  $ws.Cells.Item(55, 8).Value = 504.9565
  $ws.Cells.Item(55, 9).Value = 423
  $ws.Cells.Item(55, 10).Value = 800
  $ws.Cells.Item(55, 11).Value = 423
  $ws.Cells.Item(55, 12).Value = 800
  $ws.Cells.Item(55, 13).Value = -250
  $ws.Cells.Item(55, 14).Value = -1146
  # hunk
  $ws.Cells.Item(82, 8).Value = 1423.1
  $ws.Cells.Item(82, 9).Value = 1119.8572
  $ws.Cells.Item(82, 10).Value = 2130.6667
  $ws.Cells.Item(82, 11).Value = 1119.8572
  $ws.Cells.Item(82, 12).Value = 2130.6667
  $ws.Cells.Item(82, 13).Value = -758.8571999999999
  $ws.Cells.Item(82, 14).Value = -2852.6667
  # hunk
  $ws.Cells.Item(85, 8).Value = 1423.1
  $ws.Cells.Item(85, 9).Value = 1119.8572
  $ws.Cells.Item(85, 10).Value = 2130.6667
  $ws.Cells.Item(85, 11).Value = 1119.8572
  $ws.Cells.Item(85, 12).Value = 2130.6667
  $ws.Cells.Item(85, 13).Value = 128.1428000000001
  $ws.Cells.Item(85, 14).Value = -4626.6667
  # hunk
  $ws.Cells.Item(100, 8).Value = 2834
  $ws.Cells.Item(100, 10).Value = 3249.5
  $ws.Cells.Item(100, 12).Value = 3249.5
  $ws.Cells.Item(100, 14).Value = -4331.5
  # hunk
  $ws.Cells.Item(122, 8).Value = 2066.3333
  $ws.Cells.Item(122, 9).Value = 1849.5
  $ws.Cells.Item(122, 10).Value = 2500
  $ws.Cells.Item(122, 11).Value = 5548.5
  $ws.Cells.Item(122, 12).Value = 7500
  $ws.Cells.Item(122, 13).Value = -3098.5
  $ws.Cells.Item(122, 14).Value = -12400
  # hunk
  $ws.Cells.Item(132, 8).Value = 11596.19
  $ws.Cells.Item(132, 10).Value = 8851.333000000001
  $ws.Cells.Item(132, 12).Value = 26553.999
  $ws.Cells.Item(132, 14).Value = -31613.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
  # hunk
  $ws.Cells.Item(132, 8).Value = 4171.026
  $ws.Cells.Item(132, 9).Value = 4134.853
  $ws.Cells.Item(132, 10).Value = 4417
  $ws.Cells.Item(132, 11).Value = 12404.559
  $ws.Cells.Item(132, 12).Value = 13251
  $ws.Cells.Item(132, 13).Value = -9874.559000000001
  $ws.Cells.Item(132, 14).Value = -18311
  # hunk
  $ws.Cells.Item(136, 8).Value = 19616176
  $ws.Cells.Item(136, 9).Value = 22231444
  $ws.Cells.Item(136, 10).Value = 1673.5
  $ws.Cells.Item(136, 11).Value = 66694332
  $ws.Cells.Item(136, 12).Value = 5020.5
  $ws.Cells.Item(136, 13).Value = -66691782
  $ws.Cells.Item(136, 14).Value = -10120.5
